$wb = $excel.ActiveWorkbook

# --- Warehouse_Progress sheet: update milestone target dates ---
# Dates are stored as plain text (not real dates), so force a text
# number format before assigning to stop Excel auto-converting the
# string into a date serial, then restore the default "Normal" style
# so no extra cell formatting is introduced.
$wsProgress = $wb.Worksheets.Item("Warehouse_Progress")

$wsProgress.Range("B4").NumberFormat = "@"
$wsProgress.Range("B4").Value = "2026-01-15"
$wsProgress.Range("B4").Style = "Normal"

$wsProgress.Range("B5").NumberFormat = "@"
$wsProgress.Range("B5").Value = "2026-02-15"
$wsProgress.Range("B5").Style = "Normal"

$wsProgress.Range("B6").NumberFormat = "@"
$wsProgress.Range("B6").Value = "2026-03-15"
$wsProgress.Range("B6").Style = "Normal"

$wsProgress.Range("B7").NumberFormat = "@"
$wsProgress.Range("B7").Value = "2026-04-15"
$wsProgress.Range("B7").Style = "Normal"

# --- Settings sheet: update Warehouse Target Date ---
$wsSettings = $wb.Worksheets.Item("Settings")

$wsSettings.Range("B6").NumberFormat = "@"
$wsSettings.Range("B6").Value = "2026-04-15"
$wsSettings.Range("B6").Style = "Normal"
